$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 08:16"

# Country-row updates (name in col A, stats in cols B:H)
# Source counts were re-scraped, which both changed values and
# re-sorted a handful of rows (ties / new totals reshuffled order).
$countryUpdates = @(
    @{ Row = 23; Name = "Israel"; Vals = @(2030, 100, 58, 1967, 37, 2, 5) },
    @{ Row = 52; Name = "Mexico"; Vals = @(405, 38, 4, 396, 1, 1, 5) },
    @{ Row = 53; Name = "Egipto"; Vals = @(402, 0, 80, 302, 0, 0, 20) },
    @{ Row = 54; Name = "Barein"; Vals = @(392, 0, 177, 212, 2, 0, 3) },
    @{ Row = 55; Name = "Argentina"; Vals = @(387, 0, 52, 329, 0, 0, 6) },
    @{ Row = 56; Name = "Hong Kong"; Vals = @(387, 0, 102, 281, 4, 0, 4) },
    @{ Row = 57; Name = "Croacia"; Vals = @(382, 0, 5, 376, 6, 0, 1) },
    @{ Row = 58; Name = "Colombia"; Vals = @(378, 0, 6, 369, 0, 0, 3) },
    @{ Row = 59; Name = "Estonia"; Vals = @(369, 0, 7, 362, 4, 0, 0) },
    @{ Row = 112; Name = "Montenegro"; Vals = @(52, 5, 0, 51, 0, 0, 1) },
    @{ Row = 113; Name = "Liechtenstein"; Vals = @(51, 0, 0, 51, 0, 0, 0) },
    @{ Row = 114; Name = "Cuba"; Vals = @(48, 0, 1, 46, 2, 0, 1) },
    @{ Row = 121; Name = "Banglades"; Vals = @(39, 0, 7, 27, 0, 1, 5) },
    @{ Row = 136; Name = "Madagascar"; Vals = @(19, 2, 0, 19, 0, 0, 0) },
    @{ Row = 137; Name = "Barbados"; Vals = @(18, 0, 0, 18, 0, 0, 0) },
    @{ Row = 138; Name = "Islas Virgenes de los Estados Unidos"; Vals = @(17, 0, 0, 17, 0, 0, 0) },
    @{ Row = 180; Name = "Nepal"; Vals = @(3, 1, 1, 2, 0, 0, 0) },
    @{ Row = 181; Name = "Gambia"; Vals = @(3, 0, 0, 2, 0, 0, 1) },
    @{ Row = 182; Name = "Zimbabue"; Vals = @(3, 0, 0, 2, 0, 0, 1) },
    @{ Row = 183; Name = "Sudan"; Vals = @(3, 0, 0, 2, 0, 0, 1) },
    @{ Row = 184; Name = "Cabo Verde"; Vals = @(3, 0, 0, 2, 0, 0, 1) },
    @{ Row = 186; Name = "Butan"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 187; Name = "Nicaragua"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 188; Name = "Mauritania"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 189; Name = "San Martin (Parte Holandesa)"; Vals = @(2, 0, 0, 2, 0, 0, 0) }
)

foreach ($item in $countryUpdates) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Name
    $cols = @(2, 3, 4, 5, 6, 7, 8)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $item.Vals[$i]
    }
}
